$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Octubre de 2020 a las 08:14"

# --- Row 26/27: Ucrania overtakes Pakistan ---
$ws.Range("A26").Value = "Ucrania"
$ws.Range("B26").Value = 330396
$ws.Range("C26").Value = 7517
$ws.Range("D26").Value = 137578
$ws.Range("E26").Value = 186654
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 121
$ws.Range("H26").Value = 6164

$ws.Range("A27").Value = "Pakistan"
$ws.Range("B27").Value = 326216
$ws.Range("C27").Value = 736
$ws.Range("D27").Value = 309646
$ws.Range("E27").Value = 9855
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 13
$ws.Range("H27").Value = 6715

# --- Row 28: Israel updated counts ---
$ws.Range("B28").Value = 308572
$ws.Range("C28").Value = 325
$ws.Range("D28").Value = 288973
$ws.Range("E28").Value = 17280

# --- Row 63: Uzbekistan updated counts ---
$ws.Range("B63").Value = 64633
$ws.Range("C63").Value = 194
$ws.Range("D63").Value = 61734
$ws.Range("E63").Value = 2358
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 541

# --- Row 68/69: Kirguistan overtakes Irlanda ---
$ws.Range("A68").Value = "Kirguistan"
$ws.Range("B68").Value = 54588
$ws.Range("C68").Value = 582
$ws.Range("D68").Value = 47050
$ws.Range("E68").Value = 6412
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 1126

$ws.Range("A69").Value = "Irlanda"
$ws.Range("B69").Value = 54476
$ws.Range("C69").Value = 0
$ws.Range("D69").Value = 23364
$ws.Range("E69").Value = 29241
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 1871

# --- Row 85: El Salvador updated counts ---
$ws.Range("E85").Value = 3418
$ws.Range("G85").Value = 4
$ws.Range("H85").Value = 940

# --- Row 216/217: Islas Malvinas now listed ahead of Montserrat ---
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 13
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0

$ws.Range("A217").Value = "Montserrat"
$ws.Range("B217").Value = 13
$ws.Range("C217").Value = 0
$ws.Range("D217").Value = 12
$ws.Range("E217").Value = 0
$ws.Range("F217").Value = 0
$ws.Range("G217").Value = 0
$ws.Range("H217").Value = 1
